$d = $word.ActiveDocument

function Replace-Text($old, $new) {
    $d.Content.Find.Execute($old, $true, $false, $false, $false, $false, $true, 1, $false, $new, 2)
}

Replace-Text "2024-09-10 Tuesday" "2024-09-11 Wednesday"

Replace-Text "429÷6=71, 3" "647÷5=129, 2"
Replace-Text "890÷9=98, 8" "779÷2=389, 1"
Replace-Text "812÷3=270, 2" "278÷3=92, 2"
Replace-Text "960÷7=137, 1" "660÷5=132, 0"
Replace-Text "877÷8=109, 5" "814÷6=135, 4"

Replace-Text "301÷2=150, 1" "381÷6=63, 3"
Replace-Text "108÷9=12, 0" "141÷8=17, 5"
Replace-Text "782÷3=260, 2" "389÷5=77, 4"
Replace-Text "688÷9=76, 4" "325÷3=108, 1"
Replace-Text "581÷4=145, 1" "544÷9=60, 4"

Replace-Text "328÷7=46, 6" "959÷9=106, 5"
Replace-Text "674÷6=112, 2" "479÷8=59, 7"
Replace-Text "693÷6=115, 3" "103÷8=12, 7"
Replace-Text "449÷7=64, 1" "796÷3=265, 1"
Replace-Text "889÷7=127, 0" "238÷7=34, 0"

Replace-Text "581÷7=83, 0" "274÷8=34, 2"
Replace-Text "975÷9=108, 3" "238÷3=79, 1"
Replace-Text "509÷3=169, 2" "984÷4=246, 0"
Replace-Text "335÷2=167, 1" "821÷4=205, 1"
Replace-Text "972÷5=194, 2" "246÷5=49, 1"

Replace-Text "830÷2=415, 0" "236÷4=59, 0"
Replace-Text "830÷4=207, 2" "139÷3=46, 1"
Replace-Text "751÷4=187, 3" "958÷7=136, 6"
Replace-Text "494÷7=70, 4" "866÷3=288, 2"
Replace-Text "915÷6=152, 3" "541÷6=90, 1"

Write-Output "Done"
